# Add a "promo" and "fecha final de promo" column to the products sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column headers (F1, G1) -> new shared strings "promo" / "fecha final de promo"
$ws.Range("F1").Value = "promo"
$ws.Range("G1").Value = "fecha final de promo"

# Every product row (2-78) gets a promo value of 85
$ws.Range("F2:F78").Value = 85

# Only the first few rows (2-6) have a promo end date of "Diciembre"
$ws.Range("G2:G6").Value = "Diciembre"

# Fix mis-copied category values on rows 9-11 (were all "Piñatas")
$ws.Range("D9").Value = "Juguetes"
$ws.Range("D10").Value = "Sorpresas"
$ws.Range("D11").Value = "Dulces"

# Resize the columns to fit their (new) contents, as Excel's "best fit" would
$ws.Columns.Item(2).ColumnWidth = 4.585
$ws.Columns.Item(3).ColumnWidth = 54.585
$ws.Columns.Item(4).ColumnWidth = 7.585
$ws.Columns.Item(5).ColumnWidth = 18.752
$ws.Columns.Item(7).ColumnWidth = 18.585

# Selection left on D12 when the file was saved
[void]$ws.Range("D12").Select()
